{"js": "const replacements = [\n  [\"2024-11-03 Sunday\", \"2024-11-04 Monday\"],\n  [\"383\u00d72=\", \"579\u00d78=\"],\n  [\"530\u00d78=\", \"502\u00d73=\"],\n  [\"705\u00d72=\", \"277\u00d77=\"],\n  [\"425\u00d78=\", \"455\u00d77=\"],\n  [\"146\u00d76=\", \"222\u00d75=\"],\n  [\"296\u00d73=\", \"779\u00d74=\"],\n  [\"333\u00d79=\", \"101\u00d73=\"],\n  [\"382\u00d75=\", \"236\u00d77=\"],\n  [\"607\u00d77=\", \"398\u00d78=\"],\n  [\"661\u00d76=\", \"124\u00d79=\"],\n  [\"883\u00d72=\", \"502\u00d79=\"],\n  [\"424\u00d78=\", \"827\u00d73=\"],\n  [\"341\u00d72=\", \"693\u00d73=\"],\n  [\"326\u00d78=\", \"921\u00d75=\"],\n  [\"259\u00d79=\", \"387\u00d73=\"],\n  [\"448\u00d79=\", \"457\u00d76=\"],\n  [\"308\u00d77=\", \"790\u00d74=\"],\n  [\"503\u00d79=\", \"103\u00d72=\"],\n  [\"133\u00d77=\", \"879\u00d74=\"],\n  [\"782\u00d78=\", \"258\u00d79=\"],\n  [\"925\u00d74=\", \"841\u00d72=\"],\n  [\"220\u00d74=\", \"371\u00d75=\"],\n  [\"308\u00d79=\", \"840\u00d77=\"],\n  [\"292\u00d79=\", \"604\u00d73=\"],\n  [\"118\u00d74=\", \"190\u00d76=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-03 Sunday\", \"2024-11-04 Monday\"),\n    @(\"383\u00d72=\", \"579\u00d78=\"),\n    @(\"530\u00d78=\", \"502\u00d73=\"),\n    @(\"705\u00d72=\", \"277\u00d77=\"),\n    @(\"425\u00d78=\", \"455\u00d77=\"),\n    @(\"146\u00d76=\", \"222\u00d75=\"),\n    @(\"296\u00d73=\", \"779\u00d74=\"),\n    @(\"333\u00d79=\", \"101\u00d73=\"),\n    @(\"382\u00d75=\", \"236\u00d77=\"),\n    @(\"607\u00d77=\", \"398\u00d78=\"),\n    @(\"661\u00d76=\", \"124\u00d79=\"),\n    @(\"883\u00d72=\", \"502\u00d79=\"),\n    @(\"424\u00d78=\", \"827\u00d73=\"),\n    @(\"341\u00d72=\", \"693\u00d73=\"),\n    @(\"326\u00d78=\", \"921\u00d75=\"),\n    @(\"259\u00d79=\", \"387\u00d73=\"),\n    @(\"448\u00d79=\", \"457\u00d76=\"),\n    @(\"308\u00d77=\", \"790\u00d74=\"),\n    @(\"503\u00d79=\", \"103\u00d72=\"),\n    @(\"133\u00d77=\", \"879\u00d74=\"),\n    @(\"782\u00d78=\", \"258\u00d79=\"),\n    @(\"925\u00d74=\", \"841\u00d72=\"),\n    @(\"220\u00d74=\", \"371\u00d75=\"),\n    @(\"308\u00d79=\", \"840\u00d77=\"),\n    @(\"292\u00d79=\", \"604\u00d73=\"),\n    @(\"118\u00d74=\", \"190\u00d76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
